# ADD: new Scrum, FIX: Handleidingen
#
# Fix the "Done!" heading in the installation manual: the heading
# consists of a "Done" run (wrapped in spell-check proofErr markers)
# followed by a separate run containing just "!". Remove that trailing
# "!" run so the heading reads "Done", without disturbing the "Done"
# run or the surrounding <w:proofErr/> elements.

$d = $word.ActiveDocument

# Locate the "Done!" heading first so we only ever touch the "!" that
# immediately follows it (not some unrelated punctuation elsewhere).
$heading = $d.Content
$headingFound = $heading.Find.Execute(
    "Done!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)

if ($headingFound) {
    # Narrow the range to just the trailing "!" inside the matched text.
    $bang = $d.Range($heading.End - 1, $heading.End)

    if ($bang.Text -eq "!") {
        # Deleting the range removes the text; because it is the run's
        # entire content, the now-empty <w:r> is dropped on save instead
        # of being left behind as an empty run.
        $bang.Delete()
    }
}
